$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-7 from 45212 to 45221
$newValue = 45221
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
